$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.489.67"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("E3").Value = "  +1.75%  "
$ws.Range("E4").Value = "  -0.05%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "514.82"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.58%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "153.61"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.13%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.598"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +4.66%  "
$ws.Range("D9").Value = "2.617.33"
$ws.Range("E9").Value = "  +0.78%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "6.69"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +3.47%  "
$ws.Range("E11").Value = "  +1.14%  "
$ws.Range("E12").Value = "  +2.01%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.130"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("D14").Value = "3.061.89"
$ws.Range("E14").Value = "  +1.53%  "
$ws.Range("D15").Value = "60.566.79"
$ws.Range("E15").Value = "  +1.26%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "21.65"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("E17").Value = "  +2.05%  "
$ws.Range("D18").Value = "2.609.63"
$ws.Range("E18").Value = "  +1.06%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.75"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "357.77"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +5.97%  "
$ws.Range("E21").Value = "  +3.48%  "
$ws.Range("E22").Value = "  +2.61%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "61.01"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.32%  "
$ws.Range("E25").Value = "  +1.82%  "
$ws.Range("D26").Value = "2.725.09"
$ws.Range("E26").Value = "  +1.31%  "
$ws.Range("E27").Value = "  +1.46%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").Value = "0.0₃0840"
$ws.Range("E29").Value = "  -0.30%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "7.33"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("E31").Value = "  +0.23%  "
$ws.Range("E32").Value = "  +1.79%  "
$ws.Range("E33").Value = "  +2.55%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "5.92"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +4.22%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "150.54"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -3.55%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.01"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.920"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +8.76%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.19"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  +1.52%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.846"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.41%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "36.24"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +2.03%  "
$ws.Range("E42").Value = "  +0.19%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "289.98"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -2.02%  "
$ws.Range("E44").Value = "  +2.48%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.620"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.43%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.25%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0556"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -1.32%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "19.63"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.96"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +1.95%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.0237"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.72%  "
$ws.Range("E51").Value = "  +0.46%  "
